$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.319.23"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").Value = "2.175.83"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.68"
$ws.Range("E5").Value = "  +6.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.05"
$ws.Range("E7").Value = "  +3.79%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.01"
$ws.Range("E10").Value = "  +3.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").Value = "2.501.81"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.16"
$ws.Range("E15").Value = "  -1.27%  "

$ws.Range("D16").Value = "2.172.07"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "42.219.90"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.68"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.86"
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.93"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("E23").Value = "  +5.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.50"
$ws.Range("E24").Value = "  -5.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.44"
$ws.Range("E26").Value = "  -2.37%  "

$ws.Range("E27").Value = "  +2.20%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  -1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.74"
$ws.Range("E30").Value = "  +12.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.20"
$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0813"
$ws.Range("E33").Value = "  +5.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("E34").Value = "  -3.16%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  +3.96%  "

$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0333"
$ws.Range("E38").Value = "  +7.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.83"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.196"
$ws.Range("E41").Value = "  +3.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.44"
$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.15"
$ws.Range("E43").Value = "  -3.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.16"
$ws.Range("E44").Value = "  +6.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.470"
$ws.Range("E45").Value = "  +15.47%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0972"
$ws.Range("E46").Value = "  +0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.25"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +10.31%  "

$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  +0.38%  "
